$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 460, shifting existing rows 460..551 down to 461..552
# (mirrors Excel's native "Insert Row" which shifts cells down and
# copies formatting from the row above for the newly created row).
$ws.Rows(460).Insert()

# Populate the newly inserted row 460 with the new weekly record.
$ws.Cells.Item(460, 1).Value() = 8
$ws.Cells.Item(460, 2).Value() = "Terminal La Palmera de La Serena"
$ws.Cells.Item(460, 3).Value() = "Coquimbo"
$ws.Cells.Item(460, 4).Value() = 45209
$ws.Cells.Item(460, 5).Value() = 4
$ws.Cells.Item(460, 6).Value() = 100112003
$ws.Cells.Item(460, 7).Value() = "Ajo"
$ws.Cells.Item(460, 8).Value() = "Chino"
$ws.Cells.Item(460, 9).Value() = "Primera"
$ws.Cells.Item(460, 10).Value() = 360
$ws.Cells.Item(460, 11).Value() = 21500
$ws.Cells.Item(460, 12).Value() = 22000
$ws.Cells.Item(460, 13).Value() = 21750
$ws.Cells.Item(460, 14).Value() = "$/caja 10 kilos"
$ws.Cells.Item(460, 15).Value() = "China"
$ws.Cells.Item(460, 16).Value() = 2175
$ws.Cells.Item(460, 17).Value() = 10
$ws.Cells.Item(460, 18).Value() = "Hortaliza"
